$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $escaped = $text -replace '"', '""'
    $ws.Range($range).Formula = '="' + $escaped + '"'
    $ws.Range($range).Copy()
    $ws.Range($range).PasteSpecial(-4163)
}

Set-TextValue "D2" '30.232.07'
Set-TextValue "E2" '  +3.27%  '
Set-TextValue "D3" '1.908.76'
Set-TextValue "E3" '  +0.33%  '
Set-TextValue "E4" '  -0.38%  '
Set-TextValue "D5" '326.42'
Set-TextValue "E5" '  +3.59%  '
Set-TextValue "E6" '  -0.32%  '
Set-TextValue "D7" '0.5157'
Set-TextValue "E7" '  +0.29%  '
Set-TextValue "D8" '0.4032'
Set-TextValue "E8" '  +2.45%  '
Set-TextValue "D9" '0.08477'
Set-TextValue "E9" '  +0.02%  '
Set-TextValue "D10" '42.71'
Set-TextValue "D11" '1.118'
Set-TextValue "E11" '  -0.12%  '
Set-TextValue "D12" '23.44'
Set-TextValue "E12" '  +13.52%  '
Set-TextValue "D13" '6.459'
Set-TextValue "E13" '  +3.11%  '
Set-TextValue "D14" '1.912.03'
Set-TextValue "E14" '  +0.47%  '
Set-TextValue "D15" '7.367'
Set-TextValue "E15" '  +0.09%  '
Set-TextValue "E16" '  -0.32%  '
Set-TextValue "D17" '95.12'
Set-TextValue "E17" '  +2.01%  '
Set-TextValue "E18" '  +0.69%  '
Set-TextValue "D19" '0.06683'
Set-TextValue "E19" '  -0.73%  '
Set-TextValue "D20" '18.37'
Set-TextValue "E20" '  +2.46%  '
Set-TextValue "E22" '  -0.69%  '
Set-TextValue "D23" '30.227.11'
Set-TextValue "E23" '  +3.22%  '
Set-TextValue "D24" '11.30'
Set-TextValue "E24" '  +1.14%  '
Set-TextValue "D25" '2.222'
Set-TextValue "E25" '  +0.31%  '
Set-TextValue "D26" '2.146.09'
Set-TextValue "E26" '  +1.23%  '
Set-TextValue "D27" '21.74'
Set-TextValue "E27" '  +3.92%  '
Set-TextValue "D28" '161.66'
Set-TextValue "E28" '  +0.86%  '
Set-TextValue "E29" '  -2.49%  '
Set-TextValue "D30" '129.62'
Set-TextValue "E30" '  +1.75%  '
Set-TextValue "D31" '1.100'
Set-TextValue "E31" '  +3.84%  '
Set-TextValue "E32" '  +1.02%  '
Set-TextValue "D33" '6.061'
Set-TextValue "E33" '  +0.09%  '
Set-TextValue "D34" '3.755'
Set-TextValue "E34" '  +2.96%  '
Set-TextValue "D35" '0.02505'
Set-TextValue "E35" '  +1.00%  '
Set-TextValue "D36" '0.06592'
Set-TextValue "E36" '  -0.08%  '
Set-TextValue "D37" '0.2214'
Set-TextValue "E37" '  +0.67%  '
Set-TextValue "D38" '5.239'
Set-TextValue "E38" '  +2.16%  '
Set-TextValue "D39" '1.238'
Set-TextValue "E39" '  +0.04%  '
Set-TextValue "D40" '11.93'
Set-TextValue "E40" '  +5.69%  '
Set-TextValue "D41" '8.816'
Set-TextValue "E41" '  -3.30%  '
Set-TextValue "D42" '0.6521'
Set-TextValue "E42" '  -0.03%  '
Set-TextValue "E43" '  +0.17%  '
Set-TextValue "D44" '0.6133'
Set-TextValue "E44" '  +1.26%  '
Set-TextValue "D45" '13.21'
Set-TextValue "E45" '  -0.15%  '
Set-TextValue "D46" '3.717'
Set-TextValue "E46" '  +1.06%  '
Set-TextValue "E47" '  +0.38%  '
Set-TextValue "D48" '1.247'
Set-TextValue "E48" '  +1.31%  '
Set-TextValue "D49" '125.17'
Set-TextValue "E49" '  +1.81%  '
Set-TextValue "D50" '1.159'
Set-TextValue "E50" '  -1.51%  '
Set-TextValue "E51" '  +2.17%  '

$ws.Range("A1").Select() | Out-Null
